# Scheduled runner update: refresh market-price snapshot columns (H-N)
# for the affected Leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 3220.2222
$ws.Range("I80").Value = 522.8461
$ws.Range("J80").Value = 5724.9287
$ws.Range("K80").Value = 1568.5383
$ws.Range("L80").Value = 17174.7861
$ws.Range("M80").Value = -570.5382999999999
$ws.Range("N80").Value = -19170.7861

$ws.Range("H83").Value = 3220.2222
$ws.Range("I83").Value = 522.8461
$ws.Range("J83").Value = 5724.9287
$ws.Range("K83").Value = 4705.6149
$ws.Range("L83").Value = 51524.35830000001
$ws.Range("M83").Value = 286.3851000000004
$ws.Range("N83").Value = -61508.35830000001

$ws.Range("H86").Value = 3750.75
$ws.Range("I86").Value = 3000
$ws.Range("J86").Value = 4001
$ws.Range("K86").Value = 3000
$ws.Range("L86").Value = 4001
$ws.Range("M86").Value = -1877
$ws.Range("N86").Value = -6247

$ws.Range("H88").Value = 7294.5
$ws.Range("J88").Value = 8828.799999999999
$ws.Range("L88").Value = 8828.799999999999
$ws.Range("N88").Value = -9640.799999999999

$ws.Range("H89").Value = 3750.75
$ws.Range("I89").Value = 3000
$ws.Range("J89").Value = 4001
$ws.Range("K89").Value = 15000
$ws.Range("L89").Value = 20005
$ws.Range("M89").Value = -9384
$ws.Range("N89").Value = -31237

$ws.Range("H91").Value = 7294.5
$ws.Range("J91").Value = 8828.799999999999
$ws.Range("L91").Value = 8828.799999999999
$ws.Range("N91").Value = -11636.8

$ws.Range("H116").Value = 5360.353
$ws.Range("I116").Value = 5200
$ws.Range("J116").Value = 5540.75
$ws.Range("K116").Value = 5200
$ws.Range("L116").Value = 5540.75
$ws.Range("M116").Value = -1758
$ws.Range("N116").Value = -12424.75

$ws.Range("H132").Value = 10005120
$ws.Range("I132").Value = 12504525
$ws.Range("K132").Value = 37513575
$ws.Range("M132").Value = -37511045

$ws.Range("H137").Value = 3708486.2
$ws.Range("J137").Value = 3625
$ws.Range("L137").Value = 10875
$ws.Range("N137").Value = -15975

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6389.656
$ws.Range("I32").Value = 5163.724
$ws.Range("J32").Value = 24165.666
$ws.Range("K32").Value = 5163.724
$ws.Range("L32").Value = 24165.666
$ws.Range("M32").Value = -4876.724
$ws.Range("N32").Value = -24739.666

$ws.Range("H63").Value = 2494
$ws.Range("I63").Value = 1825.9333
$ws.Range("K63").Value = 1825.9333
$ws.Range("M63").Value = -1139.9333

$ws.Range("H66").Value = 2494
$ws.Range("I66").Value = 1825.9333
$ws.Range("K66").Value = 9129.666499999999
$ws.Range("M66").Value = -5697.666499999999

$ws.Range("H110").Value = 1630.9445
$ws.Range("I110").Value = 447.0909
$ws.Range("J110").Value = 3491.2856
$ws.Range("K110").Value = 447.0909
$ws.Range("L110").Value = 3491.2856
$ws.Range("M110").Value = 1597.9091
$ws.Range("N110").Value = -7581.2856

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1188.4814
$ws.Range("I20").Value = 1212.2222
$ws.Range("J20").Value = 1141
$ws.Range("K20").Value = 1212.2222
$ws.Range("L20").Value = 1141
$ws.Range("M20").Value = -965.2221999999999
$ws.Range("N20").Value = -1635

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2327971.5
$ws.Range("I31").Value = 2779479.5
$ws.Range("K31").Value = 2779479.5
$ws.Range("M31").Value = -2779184.5

$ws.Range("H34").Value = 2327971.5
$ws.Range("I34").Value = 2779479.5
$ws.Range("K34").Value = 2779479.5
$ws.Range("M34").Value = -2779277.5

$ws.Range("H58").Value = 35720316
$ws.Range("I58").Value = 3789.1428
$ws.Range("J58").Value = 71436840
$ws.Range("K58").Value = 3789.1428
$ws.Range("L58").Value = 71436840
$ws.Range("M58").Value = -3586.1428
$ws.Range("N58").Value = -71437246

$ws.Range("H99").Value = 2217.318
$ws.Range("I99").Value = 1833.1111
$ws.Range("J99").Value = 2483.3076
$ws.Range("K99").Value = 1833.1111
$ws.Range("L99").Value = 2483.3076
$ws.Range("M99").Value = -335.1111000000001
$ws.Range("N99").Value = -5479.3076

$ws.Range("H126").Value = 2217.318
$ws.Range("I126").Value = 1833.1111
$ws.Range("J126").Value = 2483.3076
$ws.Range("K126").Value = 5499.3333
$ws.Range("L126").Value = 7449.9228
$ws.Range("M126").Value = -3029.3333
$ws.Range("N126").Value = -12389.9228

$ws.Range("H132").Value = 2651.6
$ws.Range("I132").Value = 2633.3333
$ws.Range("K132").Value = 7899.999899999999
$ws.Range("M132").Value = -5369.999899999999

$ws.Range("H134").Value = 2854.0667
$ws.Range("I134").Value = 1179.7
$ws.Range("K134").Value = 3539.1
$ws.Range("M134").Value = -1004.1

$ws.Range("H136").Value = 35720316
$ws.Range("I136").Value = 3789.1428
$ws.Range("J136").Value = 71436840
$ws.Range("K136").Value = 11367.4284
$ws.Range("L136").Value = 214310520
$ws.Range("M136").Value = -8817.428400000001
$ws.Range("N136").Value = -214315620

$ws.Range("H141").Value = 14647.368
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 14647.368
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 14647.368
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -25007.368

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 122.8
$ws.Range("I33").Value = 100
$ws.Range("J33").Value = 128.5
$ws.Range("K33").Value = 600
$ws.Range("L33").Value = 771
$ws.Range("M33").Value = -317
$ws.Range("N33").Value = -1337

$ws.Range("H122").Value = 1396.8235
$ws.Range("I122").Value = 469.25
$ws.Range("J122").Value = 2221.3333
$ws.Range("K122").Value = 4223.25
$ws.Range("L122").Value = 19991.9997
$ws.Range("M122").Value = -1773.25
$ws.Range("N122").Value = -24891.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4969.04
$ws.Range("I132").Value = 6475.273
$ws.Range("J132").Value = 3785.5715
$ws.Range("K132").Value = 19425.819
$ws.Range("L132").Value = 11356.7145
$ws.Range("M132").Value = -16895.819
$ws.Range("N132").Value = -16416.7145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 62501524
$ws.Range("I16").Value = 90910240
$ws.Range("J16").Value = 2340
$ws.Range("K16").Value = 90910240
$ws.Range("L16").Value = 2340
$ws.Range("M16").Value = -90910070
$ws.Range("N16").Value = -2680

$ws.Range("H76").Value = 30048
$ws.Range("J76").Value = 30048
$ws.Range("L76").Value = 30048
$ws.Range("N76").Value = -30724

$ws.Range("H79").Value = 30048
$ws.Range("J79").Value = 30048
$ws.Range("L79").Value = 30048
$ws.Range("N79").Value = -32388

$ws.Range("H100").Value = 1773.9
$ws.Range("I100").Value = 1136.5
$ws.Range("K100").Value = 1136.5
$ws.Range("M100").Value = -595.5
